$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values - prefix with apostrophe to force text and avoid numeric auto-parsing
$ws.Range("D2").Value = '''24.862.88'
$ws.Range("D3").Value = '''1.708.62'
$ws.Range("D5").Value = '''312.63'
$ws.Range("D6").Value = '''0.9992'
$ws.Range("D7").Value = '''0.3754'
$ws.Range("D8").Value = '''49.56'
$ws.Range("D9").Value = '''0.3464'
$ws.Range("D10").Value = '''1.217'
$ws.Range("D11").Value = '''0.07568'
$ws.Range("D13").Value = '''21.33'
$ws.Range("D14").Value = '''6.335'
$ws.Range("D15").Value = '''7.090'
$ws.Range("D16").Value = '''1.710.17'
$ws.Range("D17").Value = '''0.00001134'
$ws.Range("D18").Value = '''0.06738'
$ws.Range("D19").Value = '''0.9992'
$ws.Range("D20").Value = '''84.87'
$ws.Range("D21").Value = '''17.36'
$ws.Range("D22").Value = '''6.408'
$ws.Range("D23").Value = '''13.21'
$ws.Range("D24").Value = '''24.878.02'
$ws.Range("D25").Value = '''2.469'
$ws.Range("D26").Value = '''2.799'
$ws.Range("D27").Value = '''20.50'
$ws.Range("D28").Value = '''150.82'
$ws.Range("D29").Value = '''132.90'
$ws.Range("D30").Value = '''1.901.37'
$ws.Range("D31").Value = '''1.251'
$ws.Range("D32").Value = '''6.889'
$ws.Range("D34").Value = '''13.92'
$ws.Range("D35").Value = '''0.08853'
$ws.Range("D36").Value = '''1.761'
$ws.Range("D37").Value = '''5.668'
$ws.Range("D38").Value = '''9.351'
$ws.Range("D39").Value = '''0.06677'
$ws.Range("D40").Value = '''0.02420'
$ws.Range("D41").Value = '''0.2245'
$ws.Range("D42").Value = '''1.281'
$ws.Range("D43").Value = '''0.6478'
$ws.Range("D44").Value = '''0.9994'
$ws.Range("D45").Value = '''13.92'
$ws.Range("D46").Value = '''0.6171'
$ws.Range("D47").Value = '''3.838'
$ws.Range("D48").Value = '''2.139'
$ws.Range("D49").Value = '''130.67'
$ws.Range("D50").Value = '''0.07315'
$ws.Range("D51").Value = '''80.52'

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("E8").Value = '  +3.25%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("E10").Value = '  +3.09%  '
$ws.Range("E11").Value = '  +3.98%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  +4.43%  '
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("E15").Value = '  +4.93%  '
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("E20").Value = '  +4.66%  '
$ws.Range("E21").Value = '  +5.41%  '
$ws.Range("E22").Value = '  +4.80%  '
$ws.Range("E23").Value = '  +9.65%  '
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("E25").Value = '  +0.86%  '
$ws.Range("E26").Value = '  +4.51%  '
$ws.Range("E27").Value = '  +5.32%  '
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("E29").Value = '  +4.66%  '
$ws.Range("E30").Value = '  +2.14%  '
$ws.Range("E31").Value = '  +27.90%  '
$ws.Range("E32").Value = '  +8.58%  '
$ws.Range("E33").Value = '  +4.60%  '
$ws.Range("E34").Value = '  +10.56%  '
$ws.Range("E35").Value = '  +4.89%  '
$ws.Range("E36").Value = '  +3.65%  '
$ws.Range("E37").Value = '  +5.44%  '
$ws.Range("E38").Value = '  +4.37%  '
$ws.Range("E39").Value = '  +2.39%  '
$ws.Range("E40").Value = '  +3.88%  '
$ws.Range("E41").Value = '  +6.00%  '
$ws.Range("E42").Value = '  +1.53%  '
$ws.Range("E43").Value = '  +4.81%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  +4.59%  '
$ws.Range("E46").Value = '  +3.49%  '
$ws.Range("E47").Value = '  +1.87%  '
$ws.Range("E48").Value = '  +5.23%  '
$ws.Range("E49").Value = '  +2.63%  '
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("E51").Value = '  +5.91%  '

# Reset style on D and E columns so no extra quotePrefix/style gets recorded
$ws.Range("D2:E51").Style = "Normal"
